# Bondtech Kit reference added
# Replace the individual Bondtech BMG part links/descriptions (rows 3-8,
# "HextrudORT_HotEnd_Nova" sub-assembly) with a single consolidated
# "BMG Internals Set for HextrudORT" kit reference.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("From Fusion 360")

$newVendor = "Included in BMG Internals Set for HextrudORT"
$newMakeBuy = "(BUY) KIT"
$newUrl = "https://www.bondtech.se/product/bmg-internals-set-for-hextrudort/"

# Rows 3-8 (the Bondtech BMG parts sub-assembly block) each get the same
# consolidated vendor text and buy/kit flag.
$rows = 3,4,5,6,7,8
foreach ($r in $rows) {
    $ws.Range("K$r").Value = $newUrl
    $ws.Range("F$r").Value = $newVendor
    $ws.Range("G$r").Value = $newMakeBuy
}

# Point each of those rows' existing hyperlink at the new kit URL, leaving
# the unrelated K12 (3dpassion.com) hyperlink untouched.
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Address -ne "https://3dpassion.com/nova") {
        $hl.Address = $newUrl
    }
}

# Match the reviewer's final view state: scrolled back to the top of the
# sheet, zoomed to 100%, with I3 selected.
[void]$ws.Range("I3").Select()
$excel.ActiveWindow.Zoom = 100
